# Auto-generated edit script applying numeric corrections to Leve profit tracker sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1006620.94
$ws.Range("I17").Value = 1150.4286
$ws.Range("J17").Value = 1548028.1
$ws.Range("K17").Value = 3451.2858
$ws.Range("L17").Value = 4644084.300000001
$ws.Range("M17").Value = -3283.2858
$ws.Range("N17").Value = -4644420.300000001
$ws.Range("H40").Value = 68960.37
$ws.Range("J40").Value = 3853
$ws.Range("L40").Value = 3853
$ws.Range("N40").Value = -4203
$ws.Range("H62").Value = 8360.056
$ws.Range("I62").Value = 7642.2144
$ws.Range("K62").Value = 7642.2144
$ws.Range("M62").Value = -7018.2144
$ws.Range("H65").Value = 8360.056
$ws.Range("I65").Value = 7642.2144
$ws.Range("K65").Value = 38211.072
$ws.Range("M65").Value = -35091.072
$ws.Range("H116").Value = 76600
$ws.Range("I116").Value = 150000
$ws.Range("J116").Value = 3200
$ws.Range("K116").Value = 150000
$ws.Range("L116").Value = 3200
$ws.Range("M116").Value = -146558
$ws.Range("N116").Value = -10084
$ws.Range("H129").Value = 1561.1428
$ws.Range("I129").Value = 779.4286
$ws.Range("J129").Value = 2342.8572
$ws.Range("K129").Value = 2338.2858
$ws.Range("L129").Value = 7028.571599999999
$ws.Range("M129").Value = 2661.7142
$ws.Range("N129").Value = -17028.5716
$ws.Range("H132").Value = 8573
$ws.Range("I132").Value = 9362
$ws.Range("K132").Value = 28086
$ws.Range("M132").Value = -25556

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3078.7666
$ws.Range("I132").Value = 2774.3044
$ws.Range("K132").Value = 8322.913199999999
$ws.Range("M132").Value = -5792.913199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 27779894
$ws.Range("I134").Value = 1604
$ws.Range("J134").Value = 83336470
$ws.Range("K134").Value = 4812
$ws.Range("L134").Value = 250009410
$ws.Range("M134").Value = -2277
$ws.Range("N134").Value = -250014480
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3607.2554
$ws.Range("I31").Value = 2002.6522
$ws.Range("K31").Value = 2002.6522
$ws.Range("M31").Value = -1707.6522
$ws.Range("H34").Value = 3607.2554
$ws.Range("I34").Value = 2002.6522
$ws.Range("K34").Value = 2002.6522
$ws.Range("M34").Value = -1800.6522
$ws.Range("H86").Value = 8500
$ws.Range("I86").Value = 9000
$ws.Range("K86").Value = 9000
$ws.Range("M86").Value = -7877
$ws.Range("H89").Value = 8500
$ws.Range("I89").Value = 9000
$ws.Range("K89").Value = 45000
$ws.Range("M89").Value = -39384
$ws.Range("H97").Value = 96371
$ws.Range("J97").Value = 96371
$ws.Range("L97").Value = 96371
$ws.Range("N97").Value = -98353

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 645.9286
$ws.Range("I34").Value = 103.44444
$ws.Range("J34").Value = 1622.4
$ws.Range("K34").Value = 310.33332
$ws.Range("L34").Value = 4867.200000000001
$ws.Range("M34").Value = -226.33332
$ws.Range("N34").Value = -5035.200000000001
$ws.Range("H38").Value = 628.1
$ws.Range("I38").Value = 75.59999999999999
$ws.Range("J38").Value = 1180.6
$ws.Range("K38").Value = 226.8
$ws.Range("L38").Value = 3541.8
$ws.Range("M38").Value = 120.2
$ws.Range("N38").Value = -4235.799999999999
$ws.Range("H39").Value = 3593
$ws.Range("J39").Value = 4889.5
$ws.Range("L39").Value = 14668.5
$ws.Range("N39").Value = -15256.5
$ws.Range("H50").Value = 2812.375
$ws.Range("I50").Value = 2740
$ws.Range("K50").Value = 8220
$ws.Range("M50").Value = -7739
$ws.Range("H53").Value = 2812.375
$ws.Range("I53").Value = 2740
$ws.Range("K53").Value = 8220
$ws.Range("M53").Value = -7739
$ws.Range("H131").Value = 1283.6364
$ws.Range("I131").Value = 889.4737
$ws.Range("J131").Value = 1818.5714
$ws.Range("K131").Value = 2668.4211
$ws.Range("L131").Value = 5455.7142
$ws.Range("M131").Value = 2371.5789
$ws.Range("N131").Value = -15535.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 60500.75
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H34").Value = 95499.5
$ws.Range("I34").Value = 90999
$ws.Range("K34").Value = 90999
$ws.Range("M34").Value = -90731
$ws.Range("H76").Value = 95499.5
$ws.Range("I76").Value = 90999
$ws.Range("K76").Value = 90999
$ws.Range("M76").Value = -90684
$ws.Range("H79").Value = 95499.5
$ws.Range("I79").Value = 90999
$ws.Range("K79").Value = 90999
$ws.Range("M79").Value = -89907
$ws.Range("H113").Value = 1319.3334
$ws.Range("J113").Value = 1243.75
$ws.Range("L113").Value = 1243.75
$ws.Range("N113").Value = -5583.75
$ws.Range("H126").Value = 4426.75
$ws.Range("I126").Value = 3496.1428
$ws.Range("J126").Value = 5729.6
$ws.Range("K126").Value = 10488.4284
$ws.Range("L126").Value = 17188.8
$ws.Range("M126").Value = -8018.428400000001
$ws.Range("N126").Value = -22128.8
$ws.Range("H132").Value = 3184.2856
$ws.Range("I132").Value = 3048.3333
$ws.Range("K132").Value = 9144.999899999999
$ws.Range("M132").Value = -6614.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7743.56
$ws.Range("I7").Value = 7436.0454
$ws.Range("K7").Value = 7436.0454
$ws.Range("M7").Value = -7324.0454
$ws.Range("H16").Value = 1651.2307
$ws.Range("I16").Value = 1651.2307
$ws.Range("K16").Value = 1651.2307
$ws.Range("M16").Value = -1481.2307
$ws.Range("H76").Value = 38999.5
$ws.Range("J76").Value = 37999.5
$ws.Range("L76").Value = 37999.5
$ws.Range("N76").Value = -38675.5
$ws.Range("H79").Value = 38999.5
$ws.Range("J79").Value = 37999.5
$ws.Range("L79").Value = 37999.5
$ws.Range("N79").Value = -40339.5
$ws.Range("H99").Value = 81509.836
$ws.Range("J99").Value = 99764.75
$ws.Range("L99").Value = 99764.75
$ws.Range("N99").Value = -105754.75
$ws.Range("H102").Value = 120000
$ws.Range("J102").Value = 120000
$ws.Range("L102").Value = 120000
$ws.Range("N102").Value = -126490
$ws.Range("H126").Value = 7743.56
$ws.Range("I126").Value = 7436.0454
$ws.Range("K126").Value = 22308.1362
$ws.Range("M126").Value = -19838.1362
$ws.Range("H132").Value = 4116.3335
$ws.Range("I132").Value = 4138.6
$ws.Range("K132").Value = 12415.8
$ws.Range("M132").Value = -9885.800000000001
$ws.Range("H136").Value = 4229
$ws.Range("I136").Value = 3517.1667
$ws.Range("K136").Value = 10551.5001
$ws.Range("M136").Value = -8001.500100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 51012.668
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H93").Value = 99333
$ws.Range("J93").Value = 99333
$ws.Range("L93").Value = 99333
$ws.Range("N93").Value = -104325
$ws.Range("H100").Value = 1772.1111
$ws.Range("I100").Value = 2176.6667
$ws.Range("K100").Value = 4353.3334
$ws.Range("M100").Value = -3812.3334
$ws.Range("H102").Value = 84994.5
$ws.Range("J102").Value = 84994.5
$ws.Range("L102").Value = 84994.5
$ws.Range("N102").Value = -91484.5
$ws.Range("H113").Value = 190.6842
$ws.Range("J113").Value = 235
$ws.Range("L113").Value = 705
$ws.Range("N113").Value = -5045
$ws.Range("H118").Value = 109499
$ws.Range("J118").Value = 109499
$ws.Range("L118").Value = 109499
$ws.Range("N118").Value = -112813
$ws.Range("H132").Value = 2556.78
$ws.Range("I132").Value = 2443
$ws.Range("J132").Value = 3580.8
$ws.Range("K132").Value = 7329
$ws.Range("L132").Value = 10742.4
$ws.Range("M132").Value = -4799
$ws.Range("N132").Value = -15802.4
$ws.Range("H140").Value = 82618.336
$ws.Range("J140").Value = 82618.336
$ws.Range("L140").Value = 82618.336
$ws.Range("N140").Value = -92978.336
